$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13, shifting existing rows 13-14 down to 14-15
$ws.Rows(13).Insert()

# Populate the new row 13 with match data
$ws.Range("A13").Value2 = "ltuWJeFb"
$ws.Range("B13").Value2 = "14/11/2024"
$ws.Range("C13").Value2 = "23:35"
$ws.Range("D13").Value2 = "MEXICO - LIGA DE EXPANSION MX"
$ws.Range("E13").Value2 = "Leones Negros"
$ws.Range("F13").Value2 = "Tapatio"
$ws.Range("G13").Value2 = 1.75
$ws.Range("H13").Value2 = 3.3
$ws.Range("I13").Value2 = 4.65
$ws.Range("J13").Value2 = 2.35
$ws.Range("K13").Value2 = 2.05
$ws.Range("L13").Value2 = 4.9
$ws.Range("M13").Value2 = 1.01
$ws.Range("N13").Value2 = 8.699999999999999
$ws.Range("O13").Value2 = 1.28
$ws.Range("P13").Value2 = 3.05
$ws.Range("Q13").Value2 = 1.88
$ws.Range("R13").Value2 = 1.82
$ws.Range("S13").Value2 = 1.42
$ws.Range("T13").Value2 = 2.47
$ws.Range("U13").Value2 = 1.72
$ws.Range("V13").Value2 = 1.88
$ws.Range("W13").Value2 = 6.9
$ws.Range("X13").Value2 = 8.5
$ws.Range("Y13").Value2 = 7.9
$ws.Range("Z13").Value2 = 14.5
$ws.Range("AA13").Value2 = 13.5
$ws.Range("AB13").Value2 = 25
$ws.Range("AC13").Value2 = 9.5
$ws.Range("AD13").Value2 = 6.5
$ws.Range("AE13").Value2 = 14
$ws.Range("AF13").Value2 = 65
$ws.Range("AG13").Value2 = 500
$ws.Range("AH13").Value2 = 12.5
$ws.Range("AI13").Value2 = 28
$ws.Range("AJ13").Value2 = 15
$ws.Range("AK13").Value2 = 90
$ws.Range("AL13").Value2 = 45
$ws.Range("AM13").Value2 = 45
$ws.Range("AN13").Value2 = 3.55
$ws.Range("AO13").Value2 = 8.75
$ws.Range("AP13").Value2 = 17.5
$ws.Range("AQ13").Value2 = 32
$ws.Range("AR13").Value2 = 60
$ws.Range("AS13").Value2 = 250
$ws.Range("AT13").Value2 = 2.45
$ws.Range("AU13").Value2 = 7.1
$ws.Range("AV13").Value2 = 65
$ws.Range("AW13").Value2 = 6.3
$ws.Range("AX13").Value2 = 27
$ws.Range("AY13").Value2 = 32
$ws.Range("AZ13").Value2 = 175
$ws.Range("BA13").Value2 = 200
$ws.Range("BB13").Value2 = 400
$ws.Range("BC13").Value2 = 51
$ws.Range("BD13").Value2 = 51

# Apply corrected odds values to row 4
$ws.Range("J4").Value2 = 4.5
$ws.Range("M4").Value2 = 1.11
$ws.Range("N4").Value2 = 6.5
$ws.Range("W4").Value2 = 8.5
$ws.Range("AH4").Value2 = 6
$ws.Range("AJ4").Value2 = 10
$ws.Range("AX4").Value2 = 13

# Apply corrected odds values to row 5
$ws.Range("G5").Value2 = 1.29
$ws.Range("I5").Value2 = 13
$ws.Range("W5").Value2 = 5
$ws.Range("Y5").Value2 = 11
$ws.Range("AH5").Value2 = 19
$ws.Range("AJ5").Value2 = 34
$ws.Range("AL5").Value2 = 101
$ws.Range("AM5").Value2 = 101

# Apply corrected odds values to row 6
$ws.Range("G6").Value2 = 1.6
$ws.Range("H6").Value2 = 3.7
$ws.Range("I6").Value2 = 6
$ws.Range("J6").Value2 = 2.2
$ws.Range("M6").Value2 = 1.07
$ws.Range("N6").Value2 = 8.5
$ws.Range("Z6").Value2 = 11
$ws.Range("AE6").Value2 = 19
$ws.Range("AH6").Value2 = 13
$ws.Range("AI6").Value2 = 29
$ws.Range("AL6").Value2 = 51
$ws.Range("AP6").Value2 = 21
$ws.Range("AQ6").Value2 = 26
